$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Aggiustare la mail OK", $true, $false, $false, $false, $false, $true, 1, $false, "OK", 2)
